$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "theta_threshold_range" row (row 5) is removed entirely. Deleting the
# row shifts "pie_threshold_range" (formerly row 6) up to become row 5 and
# also takes care of updating the sheet dimension and the shared-strings
# table (the now-unused "theta_threshold_range" string drops out).
$ws.Rows(5).Delete()

# beta_distance_range's Min value changes from 5.5 to 5.
$ws.Range("B3").Value = 5

# pie_threshold_range (now row 5) gets new Min/Max values.
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 15

# The cells that used to be row 6 carried a one-off "Times New Roman" cell
# style on column B. Re-apply the plain/default formatting used by the rest
# of the table (copied from A4) so that stray style is no longer referenced.
$ws.Range("A4").Copy()
$ws.Range("B5:C5").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# The selection saved with the sheet moves to F9.
[void]$ws.Range("F9").Select()

# Page setup switches to a portrait A4 sheet.
$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1
